$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2884.818
$ws.Range("M6").Value = -912.875
$ws.Range("I6").Value = 341.625
$ws.Range("K6").Value = 1024.875
$ws.Range("L18").Value = 48499.5
$ws.Range("H18").Value = 32433
$ws.Range("J18").Value = 48499.5
$ws.Range("N18").Value = -49067.5
$ws.Range("H33").Value = 908425.5600000001
$ws.Range("I33").Value = 1149717.1
$ws.Range("M33").Value = -1149488.1
$ws.Range("K33").Value = 1149717.1
$ws.Range("H64").Value = 4283.25
$ws.Range("M64").Value = -3685.3333
$ws.Range("I64").Value = 3933.3333
$ws.Range("K64").Value = 3933.3333
$ws.Range("I67").Value = 3933.3333
$ws.Range("M67").Value = -3075.3333
$ws.Range("K67").Value = 3933.3333
$ws.Range("H67").Value = 4283.25
$ws.Range("H94").Value = 1547.1666
$ws.Range("M94").Value = -1096.1666
$ws.Range("I94").Value = 1547.1666
$ws.Range("K94").Value = 1547.1666
$ws.Range("J112").Value = 1515.6154
$ws.Range("H112").Value = 1466.2188
$ws.Range("M112").Value = -2648.4998
$ws.Range("N112").Value = -6762.8462
$ws.Range("I112").Value = 1252.1666
$ws.Range("L112").Value = 4546.8462
$ws.Range("K112").Value = 3756.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3141.7058
$ws.Range("M2").Value = -2752.6428
$ws.Range("I2").Value = 2865.6428
$ws.Range("K2").Value = 2865.6428
$ws.Range("H5").Value = 143.1875
$ws.Range("I5").Value = 120.083336
$ws.Range("M5").Value = -8.083336000000003
$ws.Range("K5").Value = 120.083336
$ws.Range("K37").Value = 16673345
$ws.Range("H37").Value = 12513758
$ws.Range("M37").Value = -16673072
$ws.Range("I37").Value = 16673345
$ws.Range("H45").Value = 28171.75
$ws.Range("J45").Value = 18220.285
$ws.Range("M45").Value = -41726.8
$ws.Range("N45").Value = -18974.285
$ws.Range("I45").Value = 42103.8
$ws.Range("L45").Value = 18220.285
$ws.Range("K45").Value = 42103.8
$ws.Range("H63").Value = 2089.5454
$ws.Range("I63").Value = 2089.5454
$ws.Range("M63").Value = -1403.5454
$ws.Range("K63").Value = 2089.5454
$ws.Range("H66").Value = 2089.5454
$ws.Range("M66").Value = -7015.726999999999
$ws.Range("I66").Value = 2089.5454
$ws.Range("K66").Value = 10447.727
$ws.Range("H116").Value = 3141.7058
$ws.Range("M116").Value = -571.6428000000001
$ws.Range("I116").Value = 2865.6428
$ws.Range("K116").Value = 2865.6428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value = 2865.6428
$ws.Range("K3").Value = 2865.6428
$ws.Range("H3").Value = 3141.7058
$ws.Range("M3").Value = -2751.6428
$ws.Range("K4").Value = 120.083336
$ws.Range("H4").Value = 143.1875
$ws.Range("M4").Value = -5.083336000000003
$ws.Range("I4").Value = 120.083336
$ws.Range("H19").Value = 900000000
$ws.Range("J19").Value = 900000000
$ws.Range("L19").Value = 900000000
$ws.Range("N19").Value = -900000346
$ws.Range("M24").Value = -6731
$ws.Range("I24").Value = 6966
$ws.Range("K24").Value = 6966
$ws.Range("H24").Value = 6966
$ws.Range("H35").Value = 19967
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H36").Value = 11558.875
$ws.Range("J36").Value = 15814.4
$ws.Range("N36").Value = -16882.4
$ws.Range("M36").Value = -3932.3335
$ws.Range("I36").Value = 4466.3335
$ws.Range("L36").Value = 15814.4
$ws.Range("K36").Value = 4466.3335
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H82").Value = 17577
$ws.Range("M82").Value = -9771
$ws.Range("I82").Value = 10154
$ws.Range("K82").Value = 10154
$ws.Range("K85").Value = 10154
$ws.Range("H85").Value = 17577
$ws.Range("M85").Value = -8828
$ws.Range("I85").Value = 10154
$ws.Range("H94").Value = 524.4
$ws.Range("M94").Value = -76.29409999999996
$ws.Range("I94").Value = 527.2941
$ws.Range("K94").Value = 527.2941
$ws.Range("H110").Value = 79998.5
$ws.Range("J110").Value = 79998.5
$ws.Range("N110").Value = -88178.5
$ws.Range("L110").Value = 79998.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 464.35294
$ws.Range("M7").Value = -252.66666
$ws.Range("I7").Value = 365.66666
$ws.Range("K7").Value = 365.66666
$ws.Range("H50").Value = 14999
$ws.Range("J50").Value = 14999
$ws.Range("N50").Value = -16249
$ws.Range("L50").Value = 14999
$ws.Range("L51").Value = 9999
$ws.Range("H51").Value = 9999
$ws.Range("J51").Value = 9999
$ws.Range("N51").Value = -11471
$ws.Range("H59").Value = 14367.25
$ws.Range("J59").Value = 14367.25
$ws.Range("N59").Value = -16657.25
$ws.Range("L59").Value = 14367.25
$ws.Range("L60").Value = 7999
$ws.Range("H60").Value = 8014.8335
$ws.Range("J60").Value = 7999
$ws.Range("N60").Value = -9021
$ws.Range("L61").Value = 9999
$ws.Range("H61").Value = 9999
$ws.Range("J61").Value = 9999
$ws.Range("N61").Value = -10695
$ws.Range("H62").Value = 3799.8
$ws.Range("J62").Value = 3825
$ws.Range("N62").Value = -5073
$ws.Range("L62").Value = 3825
$ws.Range("N65").Value = -25365
$ws.Range("L65").Value = 19125
$ws.Range("H65").Value = 3799.8
$ws.Range("J65").Value = 3825

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I3").Value = 2000.2858
$ws.Range("K3").Value = 6000.857400000001
$ws.Range("H3").Value = 2000.2858
$ws.Range("M3").Value = -5888.857400000001
$ws.Range("H7").Value = 583.3333
$ws.Range("J7").Value = 625
$ws.Range("I7").Value = 500
$ws.Range("L7").Value = 1875
$ws.Range("K7").Value = 1500
$ws.Range("M7").Value = -1388
$ws.Range("N7").Value = -2099
$ws.Range("H11").Value = 378.91428
$ws.Range("M11").Value = -889.3870899999999
$ws.Range("J11").Value = 656.25
$ws.Range("N11").Value = -2248.75
$ws.Range("I11").Value = 343.12903
$ws.Range("L11").Value = 1968.75
$ws.Range("K11").Value = 1029.38709
$ws.Range("L92").Value = 3000
$ws.Range("K92").Value = 3300
$ws.Range("H92").Value = 1025
$ws.Range("M92").Value = -2052
$ws.Range("J92").Value = 1000
$ws.Range("N92").Value = -5496
$ws.Range("I92").Value = 1100

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L111").Value = 0
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("M113").Value = -250000330
$ws.Range("I113").Value = 250002500
$ws.Range("K113").Value = 250002500
$ws.Range("H113").Value = 100004200
$ws.Range("K132").Value = 4934.700000000001
$ws.Range("H132").Value = 1837.4166
$ws.Range("M132").Value = -2404.700000000001
$ws.Range("I132").Value = 1644.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 2192.0908
$ws.Range("J20").Value = 2192.0908
$ws.Range("N20").Value = -2644.0908
$ws.Range("L20").Value = 2192.0908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K132").Value = 39679.3125
$ws.Range("H132").Value = 11556.35
$ws.Range("M132").Value = -37149.3125
$ws.Range("I132").Value = 13226.4375
